$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6225490196078431
$ws.Range("J2").Value = 0.004901960784313725
$ws.Range("P2").Value = 0.1274509803921569
$ws.Range("S2").Value = 0.06862745098039216
$ws.Range("B3").Value = 0.007575757575757576
$ws.Range("C3").Value = 0.04545454545454546
$ws.Range("P3").Value = 0.7878787878787878
$ws.Range("S3").Value = 0.1590909090909091
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.035
$ws.Range("D6").Value = 0.01
$ws.Range("F6").Value = 0.065
$ws.Range("J6").Value = 0.25
$ws.Range("O6").Value = 0.035
$ws.Range("Q6").Value = 0.12
$ws.Range("R6").Value = 0.07000000000000001
$ws.Range("S6").Value = 0.415
$ws.Range("B7").Value = 0.07246376811594203
$ws.Range("D7").Value = 0.01449275362318841
$ws.Range("F7").Value = 0.08695652173913043
$ws.Range("J7").Value = 0.08695652173913043
$ws.Range("O7").Value = 0.007246376811594203
$ws.Range("Q7").Value = 0.1376811594202899
$ws.Range("R7").Value = 0.1014492753623188
$ws.Range("S7").Value = 0.4927536231884058
$ws.Range("B8").Value = 0.06696428571428571
$ws.Range("D8").Value = 0.01785714285714286
$ws.Range("F8").Value = 0.078125
$ws.Range("J8").Value = 0.1138392857142857
$ws.Range("O8").Value = 0.01785714285714286
$ws.Range("Q8").Value = 0.1875
$ws.Range("R8").Value = 0.08482142857142858
$ws.Range("S8").Value = 0.4330357142857143
$ws.Range("B9").Value = 0.07009345794392523
$ws.Range("D9").Value = 0.004672897196261682
$ws.Range("F9").Value = 0.07476635514018691
$ws.Range("J9").Value = 0.1542056074766355
$ws.Range("O9").Value = 0.01401869158878505
$ws.Range("Q9").Value = 0.1261682242990654
$ws.Range("R9").Value = 0.1448598130841121
$ws.Range("S9").Value = 0.411214953271028
$ws.Range("B10").Value = 0.09924385633270322
$ws.Range("D10").Value = 0.01417769376181474
$ws.Range("E10").Value = 0.001890359168241966
$ws.Range("F10").Value = 0.0661625708884688
$ws.Range("J10").Value = 0.113421550094518
$ws.Range("O10").Value = 0.01890359168241966
$ws.Range("Q10").Value = 0.1852551984877127
$ws.Range("R10").Value = 0.08695652173913043
$ws.Range("S10").Value = 0.4139886578449906
$ws.Range("F11").Value = 0.004587155963302753
$ws.Range("G11").Value = 0.1559633027522936
$ws.Range("J11").Value = 0.0963302752293578
$ws.Range("K11").Value = 0.2201834862385321
$ws.Range("L11").Value = 0.518348623853211
$ws.Range("S11").Value = 0.004587155963302753
$ws.Range("G12").Value = 0.7627118644067796
$ws.Range("J12").Value = 0.2033898305084746
$ws.Range("L12").Value = 0.03389830508474576
$ws.Range("G13").Value = 0.6785714285714286
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.07142857142857142
$ws.Range("F15").Value = 0.02857142857142857
$ws.Range("H15").Value = 0.1657142857142857
$ws.Range("I15").Value = 0.1428571428571428
$ws.Range("J15").Value = 0.2342857142857143
$ws.Range("K15").Value = 0.06857142857142857
$ws.Range("M15").Value = 0.005714285714285714
$ws.Range("O15").Value = 0.04
$ws.Range("S15").Value = 0.3142857142857143
$ws.Range("F16").Value = 0.01379310344827586
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = 0.1241379310344828
$ws.Range("J16").Value = 0.4413793103448276
$ws.Range("K16").Value = 0.06896551724137931
$ws.Range("M16").Value = 0.006896551724137931
$ws.Range("O16").Value = 0.04137931034482759
$ws.Range("S16").Value = 0.103448275862069
$ws.Range("F17").Value = 0.01136363636363636
$ws.Range("H17").Value = 0.2244318181818182
$ws.Range("I17").Value = 0.1022727272727273
$ws.Range("J17").Value = 0.3892045454545455
$ws.Range("K17").Value = 0.07102272727272728
$ws.Range("M17").Value = 0.02272727272727273
$ws.Range("O17").Value = 0.0625
$ws.Range("S17").Value = 0.1164772727272727
$ws.Range("F18").Value = 0.02645502645502645
$ws.Range("H18").Value = 0.2222222222222222
$ws.Range("I18").Value = 0.07936507936507936
$ws.Range("J18").Value = 0.4126984126984127
$ws.Range("K18").Value = 0.08465608465608465
$ws.Range("M18").Value = 0.03174603174603174
$ws.Range("O18").Value = 0.0582010582010582
$ws.Range("S18").Value = 0.08465608465608465
$ws.Range("F19").Value = 0.01222707423580786
$ws.Range("H19").Value = 0.2358078602620087
$ws.Range("I19").Value = 0.1048034934497817
$ws.Range("J19").Value = 0.3851528384279476
$ws.Range("K19").Value = 0.09082969432314411
$ws.Range("M19").Value = 0.01222707423580786
$ws.Range("O19").Value = 0.05764192139737991
$ws.Range("S19").Value = 0.1013100436681223
